$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.097780108451843
$ws.Range("B1").Value = 2.486933469772339
$ws.Range("C1").Value = 6.329222679138184
$ws.Range("D1").Value = 2.20048999786377
$ws.Range("E1").Value = 1.267274856567383
